# Pending Tasks.xlsx update
# - refresh the list of pending tasks (remove/add items, reorder per latest status sort)
# - extend the tracked range through row 14 and move the Total/array-formula row to 15
# - resize Table1 to match the new data extent
# - apply the "60% - Accent2" highlight style to the Total row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old "Total" row (row 13) first -----------------------------
# The old row 13 holds the anchor of the existing CSE array formula; deleting
# the row clears that array-formula registration so a new array formula can
# be written later at its new location (row 15) without being silently
# ignored by the engine.
$ws.Rows.Item(13).Delete()

# --- Rewrite the task list (rows 2-13) --------------------------------------
$tasks = @(
    @{ Row = 2;  Name = "Fix Name in dahsboard";          Pct = 1 },
    @{ Row = 3;  Name = "Profile for each Role";           Pct = 1 },
    @{ Row = 4;  Name = "tab to navigate in dashboard";    Pct = 0 },
    @{ Row = 5;  Name = "Filter Employees in tasks";       Pct = 0 },
    @{ Row = 6;  Name = "Employee Task Tab";                Pct = 1 },
    @{ Row = 7;  Name = "Employee Home Tab";                Pct = 1 },
    @{ Row = 8;  Name = "Logout";                           Pct = 1 },
    @{ Row = 9;  Name = "Email function";                   Pct = 1 },
    @{ Row = 10; Name = "Manager Task edit";                Pct = 1 },
    @{ Row = 11; Name = "Manager Task Delete";              Pct = 1 },
    @{ Row = 12; Name = "select employee from ddl";         Pct = 1 },
    @{ Row = 13; Name = "Filter Tasks";                     Pct = 1 }
)

foreach ($t in $tasks) {
    $r = $t.Row
    $ws.Range("A$r").Value = $t.Name
    $ws.Range("B$r").Value = $t.Pct
}

# Row 13 is a freshly-created row (the old row 13 was deleted above), so it
# needs the Percent number format re-applied to match the rest of the column.
$ws.Range("B13").NumberFormat = "0%"

# --- Blank trailing row (row 14) matches the Percent formatting ------------
$ws.Range("B14").NumberFormat = "0%"

# --- Total row (row 15) ------------------------------------------------------
$ws.Range("A15").Value = "Total"
$ws.Range("B15").FormulaArray = "=SUM(B2:B14 / COUNT(B2:B14))"

# Highlight the Total row with the built-in "60% - Accent2" cell style.
$ws.Range("A15:B15").Style = "60% - Accent2"
$ws.Range("B15").NumberFormat = "0%"

$excel.CalculateFull()

# --- Resize Table1 to cover the new data extent -----------------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:B14"))

# --- Selection matches the authored workbook --------------------------------
$ws.Range("B4").Select() | Out-Null
